# Generate Report for Handoff
#
# This applies the "new handoff" report regeneration:
#  - On the Overview sheet, the row that used to describe 73dd5e31...md
#    now comes first (row 2) and 289b639f...md comes second (row 3).
#  - 289b639f...md's status moves from "Handed back: in sync with en-US"
#    to "Ready for handoff" with an updated generation timestamp.
#  - Same swap + status update happens on the per-locale sheets (zh-cn, de-de),
#    where 289b639f's row additionally gets an updated "Latest Handoff Datetime"
#    and a new stale-handback Error Detail message.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value2 = "73dd5e31-0edc-440d-85b9-4f597bbd6f29.md"
$ov.Range("B2").Value2 = "e2e\73dd5e31-0edc-440d-85b9-4f597bbd6f29.md"

$ov.Range("A3").Value2 = "289b639f-c46c-420d-b500-f629bbb747f1.md"
$ov.Range("B3").Value2 = "e2e\289b639f-c46c-420d-b500-f629bbb747f1.md"
$ov.Range("E3").Value2 = "Ready for handoff"
$ov.Range("F3").Value2 = "Ready for handoff"
$ov.Range("G3").Value2 = "2016-08-12 20:59:04"

# Recreate the hyperlinks so the underlying relationships (rId2 -> 289b639f,
# rId3 -> 73dd5e31) stay the same, only which cell/display text uses them swaps.
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/c8eb83364ce31ea55ddadd7c13f3450f1807ad7b/e2e/289b639f-c46c-420d-b500-f629bbb747f1.md", "", "", "e2e\73dd5e31-0edc-440d-85b9-4f597bbd6f29.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/c8eb83364ce31ea55ddadd7c13f3450f1807ad7b/e2e/73dd5e31-0edc-440d-85b9-4f597bbd6f29.md", "", "", "e2e\289b639f-c46c-420d-b500-f629bbb747f1.md") | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value2 = "73dd5e31-0edc-440d-85b9-4f597bbd6f29.md"
$zh.Range("G2").Value2 = "73dd5e31-0edc-440d-85b9-4f597bbd6f29.9ff5cfa3021518a1b42f6e8892800a5573b04f96.zh-cn.xlf"
$zh.Range("I2").Value2 = "73dd5e31-0edc-440d-85b9-4f597bbd6f29.md"
$zh.Range("J2").Value2 = "73dd5e31-0edc-440d-85b9-4f597bbd6f29.9ff5cfa3021518a1b42f6e8892800a5573b04f96.zh-cn.xlf"

$zh.Range("A3").Value2 = "289b639f-c46c-420d-b500-f629bbb747f1.md"
$zh.Range("C3").Value2 = "Ready for handoff"
$zh.Range("G3").Value2 = "289b639f-c46c-420d-b500-f629bbb747f1.bb4df6a04e932e392ccbb794f8ee62715002b3b9.zh-cn.xlf"
$zh.Range("H3").Value2 = "2016-08-12 20:58:54"
$zh.Range("I3").Value2 = "289b639f-c46c-420d-b500-f629bbb747f1.md"
$zh.Range("J3").Value2 = "289b639f-c46c-420d-b500-f629bbb747f1.bb4df6a04e932e392ccbb794f8ee62715002b3b9.zh-cn.xlf"
$zh.Range("P3").Value2 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/c8eb83364ce31ea55ddadd7c13f3450f1807ad7b/e2e/289b639f-c46c-420d-b500-f629bbb747f1.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/8bf6e4d5dc62c59b412067b562c08becd422a222/e2e/289b639f-c46c-420d-b500-f629bbb747f1.md."

$zh.Columns.Item(16).ColumnWidth = 39.17

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/c8eb83364ce31ea55ddadd7c13f3450f1807ad7b/e2e/289b639f-c46c-420d-b500-f629bbb747f1.md", "", "", "73dd5e31-0edc-440d-85b9-4f597bbd6f29.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/8c201d0726974504e4d68677b44b9c16a87e5508/e2e/289b639f-c46c-420d-b500-f629bbb747f1.md", "", "", "73dd5e31-0edc-440d-85b9-4f597bbd6f29.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/c8eb83364ce31ea55ddadd7c13f3450f1807ad7b/e2e/73dd5e31-0edc-440d-85b9-4f597bbd6f29.md", "", "", "289b639f-c46c-420d-b500-f629bbb747f1.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/8c201d0726974504e4d68677b44b9c16a87e5508/e2e/73dd5e31-0edc-440d-85b9-4f597bbd6f29.md", "", "", "289b639f-c46c-420d-b500-f629bbb747f1.md") | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value2 = "73dd5e31-0edc-440d-85b9-4f597bbd6f29.md"
$de.Range("G2").Value2 = "73dd5e31-0edc-440d-85b9-4f597bbd6f29.9ff5cfa3021518a1b42f6e8892800a5573b04f96.de-de.xlf"
$de.Range("I2").Value2 = "73dd5e31-0edc-440d-85b9-4f597bbd6f29.md"
$de.Range("J2").Value2 = "73dd5e31-0edc-440d-85b9-4f597bbd6f29.9ff5cfa3021518a1b42f6e8892800a5573b04f96.de-de.xlf"

$de.Range("A3").Value2 = "289b639f-c46c-420d-b500-f629bbb747f1.md"
$de.Range("C3").Value2 = "Ready for handoff"
$de.Range("G3").Value2 = "289b639f-c46c-420d-b500-f629bbb747f1.bb4df6a04e932e392ccbb794f8ee62715002b3b9.de-de.xlf"
$de.Range("H3").Value2 = "2016-08-12 20:59:04"
$de.Range("I3").Value2 = "289b639f-c46c-420d-b500-f629bbb747f1.md"
$de.Range("J3").Value2 = "289b639f-c46c-420d-b500-f629bbb747f1.bb4df6a04e932e392ccbb794f8ee62715002b3b9.de-de.xlf"
$de.Range("P3").Value2 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/c8eb83364ce31ea55ddadd7c13f3450f1807ad7b/e2e/289b639f-c46c-420d-b500-f629bbb747f1.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/8bf6e4d5dc62c59b412067b562c08becd422a222/e2e/289b639f-c46c-420d-b500-f629bbb747f1.md."

$de.Columns.Item(16).ColumnWidth = 39.17

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/c8eb83364ce31ea55ddadd7c13f3450f1807ad7b/e2e/289b639f-c46c-420d-b500-f629bbb747f1.md", "", "", "73dd5e31-0edc-440d-85b9-4f597bbd6f29.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f08befee845410e3e1ab1288b57d91a565356ce0/e2e/289b639f-c46c-420d-b500-f629bbb747f1.md", "", "", "73dd5e31-0edc-440d-85b9-4f597bbd6f29.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/c8eb83364ce31ea55ddadd7c13f3450f1807ad7b/e2e/73dd5e31-0edc-440d-85b9-4f597bbd6f29.md", "", "", "289b639f-c46c-420d-b500-f629bbb747f1.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f08befee845410e3e1ab1288b57d91a565356ce0/e2e/73dd5e31-0edc-440d-85b9-4f597bbd6f29.md", "", "", "289b639f-c46c-420d-b500-f629bbb747f1.md") | Out-Null

Write-Output "edit complete"
